# Updated boms for Word compatibility, updating user guide

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hide gridlines on the sheet view ---
$excel.ActiveWindow.DisplayGridlines = $false

# --- Replace the plain "[Add OSHPark Link]" text reference (row 3) stays text, no hyperlink ---
# (value itself unchanged, just keeping as literal text so shared string survives)

# --- Turn the D4:D11 part-link cells into real =HYPERLINK() formulas ---
$ws.Range("D4").Formula = '=HYPERLINK("https://www.amazon.com/HiLetgo-Stepstick-Stepper-Printer-Compatible/dp/B00LOF1CA2/ref=sr_1_24?crid=1H73ID64FI88C&dchild=1&keywords=a4988+stepper+motor+driver&qid=1617929080&sprefix=a4988%2Caps%2C201&sr=8-24", "Motor Driver")'
$ws.Range("D5").Formula = '=HYPERLINK("https://www.digikey.com/en/products/detail/w%C3%BCrth-elektronik/691214110002S/11477397", "2x 3.5mm Terminal")'
$ws.Range("D6").Formula = '=HYPERLINK("https://www.digikey.com/en/products/detail/w%C3%BCrth-elektronik/691214110003S/11477432","3x 3.5mm Terminal")'
$ws.Range("D7").Formula = '=HYPERLINK("https://www.digikey.com/en/products/detail/w%C3%BCrth-elektronik/860240572001/5729254","10u TH Cap")'
$ws.Range("D8").Formula = '=HYPERLINK("https://www.digikey.com/en/products/detail/stackpole-electronics-inc/CF14JT10K0/1741265","10k TH Res")'
$ws.Range("D9").Formula = '=HYPERLINK("https://www.sparkfun.com/products/16581","Pi Header")'
$ws.Range("D10").Formula = '=HYPERLINK("https://www.digikey.com/en/products/detail/harwin-inc/M20-9750446/3727931","4 Pin Right Angle")'
$ws.Range("D11").Formula = '=HYPERLINK("https://www.digikey.com/en/products/detail/sullins-connector-solutions/PPTC081LFBN-RC/810147","2x Female Headers")'

# Keep the legacy hyperlink relationship on D4 but also surface the display text (url) attribute
$ws.Range("D4").Hyperlinks.Item(1).TextToDisplay = "Motor Driver"

# --- New font: Times New Roman, used across the whole table body ---
$ws.Range("A1:D11").Font.Name = "Times New Roman"

# --- Title row: centered ---
$ws.Range("A1:D1").HorizontalAlignment = -4108

# --- Hyperlink-style cells (D4:D11) use the built-in Hyperlink cell style ---
$ws.Range("D4:D11").Style = "Hyperlink"

# --- Select the whole table as the final selection ---
$ws.Range("A1:D11").Select()
